$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - prefixed with an apostrophe so Excel keeps
# numeric-looking text (e.g. "27.00", "9.60") as a literal string instead
# of coercing it to a number; Style is then reset to Normal so no stray
# number-format/quote-prefix styling is left behind.
$priceUpdates = @{
    2 = '69.816.83'
    3 = '3.824.54'
    5 = '612.95'
    6 = '176.98'
    7 = '3.822.05'
    9 = '0.528'
    12 = '0.483'
    13 = '39.78'
    14 = '0.0000254'
    15 = '4.461.43'
    16 = '3.830.51'
    17 = '69.869.55'
    18 = '7.55'
    19 = '0.118'
    20 = '16.65'
    21 = '507.32'
    22 = '9.60'
    23 = '0.742'
    25 = '86.34'
    27 = '12.66'
    30 = '2.54'
    32 = '7.99'
    33 = '31.85'
    37 = '6.12'
    39 = '484.67'
    41 = '3.04'
    43 = '49.76'
    44 = '43.88'
    45 = '8.56'
    46 = '2.925.80'
    48 = '139.42'
    50 = '27.00'
    51 = '2.43'
}

# Volume(1h) (column E) updates - already non-numeric text (percent sign
# and padding spaces), so a plain assignment keeps them as text.
$volumeUpdates = @{
    2 = '  -0.33%  '
    3 = '  +2.21%  '
    4 = '  +0.05%  '
    5 = '  -1.56%  '
    6 = '  -1.87%  '
    7 = '  +2.30%  '
    8 = '  +0.02%  '
    9 = '  -1.18%  '
    10 = '  -0.09%  '
    11 = '  +2.56%  '
    12 = '  -0.59%  '
    13 = '  -2.96%  '
    14 = '  -1.87%  '
    15 = '  +2.14%  '
    16 = '  +2.26%  '
    17 = '  -0.25%  '
    18 = '  -0.24%  '
    19 = '  -3.19%  '
    20 = '  -0.94%  '
    21 = '  +0.20%  '
    22 = '  +2.38%  '
    23 = '  +2.52%  '
    24 = '  -2.68%  '
    25 = '  -0.47%  '
    26 = '  +4.79%  '
    27 = '  -3.41%  '
    28 = '  -6.23%  '
    29 = '  +0.04%  '
    30 = '  +2.14%  '
    31 = '  +1.64%  '
    32 = '  +0.88%  '
    33 = '  +1.95%  '
    34 = '  -1.46%  '
    35 = '  +0.04%  '
    36 = '  -1.37%  '
    37 = '  -1.45%  '
    38 = '  +5.29%  '
    39 = '  +14.21%  '
    40 = '  -0.06%  '
    41 = '  +7.17%  '
    42 = '  -2.62%  '
    43 = '  -1.40%  '
    44 = '  -2.68%  '
    45 = '  -1.72%  '
    46 = '  -2.41%  '
    47 = '  -1.10%  '
    48 = '  +0.90%  '
    49 = '  +0.05%  '
    50 = '  -1.18%  '
    51 = '  -3.18%  '
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = "'" + $priceUpdates[$row]
    $cell.Style = "Normal"
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
